# Implement HOLOVIEWSMAPPING flag for holoviews vs geoviews mapping:
# Split longitude/latitude min/max columns into separate
# "londeg"/"latdeg" (copy of longitude/latitude) plus renamed
# "londeg_min"/"londeg_max"/"latdeg_min"/"latdeg_max" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at J:K. This shifts the old
# J,K,L,M,N,O (longitude_min, longitude_max, latitude_min,
# latitude_max, date_time_min, date_time_max) columns to
# L,M,N,O,P,Q respectively, preserving their values/styles.
$ws.Range("J1:K1").EntireColumn.Insert()

# Rename the headers that shifted into their new meanings.
$ws.Range("L1").Value = "londeg_min"
$ws.Range("M1").Value = "londeg_max"
$ws.Range("N1").Value = "latdeg_min"
$ws.Range("O1").Value = "latdeg_max"

# New header labels for the freshly inserted columns.
$ws.Range("J1").Value = "londeg"
$ws.Range("K1").Value = "latdeg"

# Find the last used data row (dimension is A1:Q35 -> 35 rows).
$lastRow = $ws.UsedRange.Rows.Count

# Populate the new J (londeg) and K (latdeg) columns as copies
# of the existing G (longitude) and H (latitude) columns for
# every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($r, 8).Value()
}
